# overhaul of figure order
# Update the "Visit" column on the first worksheet ("MxM - General") so that
# the numeric visit number (2 / 3) is replaced by the text labels "V2" / "V3",
# and lower-case the column header from "Visit" to "visit".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header: B1 "Visit" -> "visit"
$ws.Cells.Item(1, 2).Value = "visit"

# Data rows 2-63 live in column B ("Visit"): turn the numeric visit id into
# its text label (2 -> "V2", 3 -> "V3").
$lastRow = 63
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -eq 2) {
        $cell.Value = "V2"
    } elseif ($val -eq 3) {
        $cell.Value = "V3"
    } elseif ($val -eq "2") {
        $cell.Value = "V2"
    } elseif ($val -eq "3") {
        $cell.Value = "V3"
    }
}
